$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the crypto price/volume table (columns D = Price, E = Volume(1h))
# with the latest GitHub Actions refresh. A couple of rows also swap which
# coin occupies them (B = Coin name, C = Link).
#
# Note: a handful of Price values are plain single-decimal numbers
# (e.g. "0.613", "235.56"). Excel's COM layer auto-coerces such strings to
# numeric cells, which would change both the cell type and the displayed
# text (e.g. "0.613" -> 0.61299999999999999). To keep these as text, exactly
# like the other Price cells in this sheet, a leading apostrophe is used to
# force text entry - Excel strips the apostrophe itself and stores the
# clean string.

$ws.Range("D2").Value = '37.311.93'
$ws.Range("E2").Value = '  +2.43%  '

$ws.Range("D3").Value = '2.059.69'
$ws.Range("E3").Value = '  +3.76%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '''235.56'
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("D6").Value = '''0.613'
$ws.Range("E6").Value = '  +2.17%  '

$ws.Range("D7").Value = '''58.05'
$ws.Range("E7").Value = '  +6.92%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  +2.83%  '

$ws.Range("D10").Value = '''57.97'
$ws.Range("E10").Value = '  -2.00%  '

$ws.Range("E11").Value = '  +2.42%  '

$ws.Range("E12").Value = '  +3.37%  '

$ws.Range("D13").Value = '2.363.52'
$ws.Range("E13").Value = '  +3.81%  '

$ws.Range("E14").Value = '  +2.09%  '

$ws.Range("E15").Value = '  +5.26%  '

$ws.Range("D16").Value = '''0.777'
$ws.Range("E16").Value = '  +3.12%  '

$ws.Range("D17").Value = '''5.19'
$ws.Range("E17").Value = '  +3.28%  '

$ws.Range("D18").Value = '2.059.89'
$ws.Range("E18").Value = '  +3.77%  '

$ws.Range("D19").Value = '37.383.01'
$ws.Range("E19").Value = '  +2.79%  '

$ws.Range("D20").Value = '''6.22'
$ws.Range("E20").Value = '  +17.79%  '

$ws.Range("E21").Value = '  +2.13%  '

$ws.Range("E22").Value = '  +1.74%  '

$ws.Range("D23").Value = '''225.97'
$ws.Range("E23").Value = '  +1.90%  '

$ws.Range("E24").Value = '  -0.03%  '

$ws.Range("D25").Value = '''2.44'
$ws.Range("E25").Value = '  +2.31%  '

$ws.Range("D26").Value = '''2.38'
$ws.Range("E26").Value = '  +0.36%  '

$ws.Range("D27").Value = '''164.43'
$ws.Range("E27").Value = '  +1.88%  '

$ws.Range("D28").Value = '''1.49'
$ws.Range("E28").Value = '  +12.85%  '

$ws.Range("D29").Value = '''8.87'
$ws.Range("E29").Value = '  +3.88%  '

$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '''0.128'
$ws.Range("E30").Value = '  +2.15%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = '''19.21'
$ws.Range("E31").Value = '  +2.30%  '

$ws.Range("E32").Value = '  +1.70%  '

$ws.Range("D33").Value = '''4.49'
$ws.Range("E33").Value = '  +3.32%  '

$ws.Range("D34").Value = '''0.0622'
$ws.Range("E34").Value = '  +2.91%  '

$ws.Range("D35").Value = '''2.56'
$ws.Range("E35").Value = '  +10.14%  '

$ws.Range("E36").Value = '  +5.73%  '

$ws.Range("E37").Value = '  +7.04%  '

$ws.Range("E38").Value = '  +0.12%  '

$ws.Range("E39").Value = '  -0.14%  '

$ws.Range("E40").Value = '  +9.27%  '

$ws.Range("B41").Value = 'Cronos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D41").Value = '''0.0985'
$ws.Range("E41").Value = '  +8.30%  '

$ws.Range("B42").Value = 'FTXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D42").Value = '''4.54'
$ws.Range("E42").Value = '  +27.81%  '

$ws.Range("E43").Value = '  -1.90%  '

$ws.Range("D44").Value = '1.474.54'

$ws.Range("D45").Value = '''96.92'
$ws.Range("E45").Value = '  +9.77%  '

$ws.Range("E46").Value = '  +6.50%  '

$ws.Range("E47").Value = '  +4.85%  '

$ws.Range("D48").Value = '''15.91'
$ws.Range("E48").Value = '  +7.70%  '

$ws.Range("E49").Value = '  +3.24%  '

$ws.Range("D50").Value = '''7.21'
$ws.Range("E50").Value = '  +6.74%  '

$ws.Range("E51").Value = '  +1.87%  '
